$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 (2-1 / Hangar 1)
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "adawd"
$ws.Range("D3").Value = "wadawd"
$ws.Range("F3").Value = "https://powerbi.bellflight.com/reports/powerbi/Piney%20Flats/Aircraft%20Services/Part%20Visibility%20Report?rs:embed=true"

# Update row 4 URLs (second link swapped for a new metrics link)
$ws.Range("F4").Value = "https://powerbi.bellflight.com/reports/powerbi/Piney%20Flats/Aircraft%20Services/Part%20Visibility%20Report?rs:embed=true|https://example.com/metrics3"

# Add new row 6 for bay 2-2
$ws.Range("A6").Value = "2-2"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "gsfsfd"
$ws.Range("D6").Value = "sfsef"
$ws.Range("E6").Value = 1
$ws.Range("F6").Formula = '=""'
